$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row
$ws.Range("A1").Value = "lang_code"
$ws.Range("B1").Value = "hierarchy_level"
$ws.Range("C1").Value = "hierarchy_level_name"
$ws.Range("D1").Value = "is_active"

# New data rows: lang_code, hierarchy_level, hierarchy_level_name, is_active
$data = @(
    @("eng", 0, "Country", $true),
    @("fra", 0, "Pays", $true),
    @("eng", 1, "Region", $true),
    @("fra", 1, "Région", $true),
    @("eng", 2, "Province", $true),
    @("fra", 2, "Province", $true),
    @("eng", 3, "City", $true),
    @("fra", 3, "Ville", $true),
    @("eng", 4, "Zone", $true),
    @("fra", 4, "Zone", $true),
    @("eng", 5, "Postal Code", $true),
    @("fra", 5, "code postal", $true)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}

# Column A data cells (A2:A13) pick up the same style as the header cells
# (bold / bordered / centered) - copy format from A1 instead of creating a
# brand-new style entry.
$ws.Range("A1").Copy()
$ws.Range("A2:A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
